$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 8570
    $ws.Range("F10").Value = 128
    $ws.Range("F13").Value = 1273
    $ws.Range("F14").Value = 284
    $ws.Range("F17").Value = 104
}
